$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 3 with the new values: text columns become a single space,
# numeric columns (I, J, K) get the new amounts.
$ws.Range("A3").Value = " "
$ws.Range("B3").Value = " "
$ws.Range("C3").Value = " "
$ws.Range("D3").Value = " "
$ws.Range("E3").Value = " "
$ws.Range("F3").Value = " "
$ws.Range("G3").Value = " "
$ws.Range("H3").Value = " "
$ws.Range("I3").Value = 24000
$ws.Range("J3").Value = 1200
$ws.Range("K3").Value = 22800

# Remove rows 4 through 9 entirely (they are no longer part of the data).
$ws.Range("A4:K9").EntireRow.Delete()
